$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the look of the
# existing header row (bold, bordered, centered) by copying the format
# from H1 ("IP") via copy/paste-special so the shared style is reused
# instead of a brand-new style entry being created.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New data values for row 2
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4
